# Edit the "Indicators" sheet of the workbook to match the author's changes:
#   - Clean up the short_name text of the "Long-term care" indicator row (E28)
#   - Toggle the Country_Profile checkboxes (column B) for a few indicators,
#     which drives the existing AutoFilter (column B == TRUE) that hides/shows rows
#   - Backfill the row-id numbers in column A for rows 162-180 (previously blank)
#   - Update the remembered selection on the Indicators sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indicators")

# 1. Fix the short_name (column E) of row 28 - replace the value that had a
#    stray non-breaking space baked into the shared string with a clean one.
$ws.Cells.Item(28, 5).Value = "Long-term care (health) expenditure"

# 2. Flip the Country_Profile checkboxes (column B) for a few indicators.
$ws.Cells.Item(31, 2).Value = $true
$ws.Cells.Item(34, 2).Value = $true
$ws.Cells.Item(44, 2).Value = $false

# 3. Re-apply the existing AutoFilter (column B, criteria TRUE) over the used
#    range so that row-hidden state is recomputed from the new checkbox values.
$filterRange = $ws.Range("B1:H180")
[void]$filterRange.AutoFilter(1, "True")

# 4. Backfill the row id numbers in column A for rows 162-180 (id = row - 1),
#    matching the pattern used by every other row on the sheet.
for ($r = 162; $r -le 180; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# 5. Update the active selection on the Indicators sheet.
[void]$ws.Range("B50").Select()
